$wb = $excel.ActiveWorkbook

$uuid1 = "19fc4d9c-7842-472d-b162-b4883640e2fb"
$uuid2 = "ffff956a9504-a356-4e11-ba86-4ab4019c890e"
$xlf1zh = "$uuid1.7d6d3d747365263bee9394c5b522fb6e1a3f6b7f.zh-cn.xlf"
$xlf1de = "$uuid1.7d6d3d747365263bee9394c5b522fb6e1a3f6b7f.de-de.xlf"
$statusText = "Ready for handoff"
$handoffDate = "2016-03-22 00:59:32"
$zhHandoffDatetime = "2016-03-22 00:59:27"
$deHandoffDatetime = "2016-03-22 00:59:32"
$neverBack = "0001-01-01 00:00:00"

# ---------- Overview sheet ----------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$uuid1.md"
$ws.Range("B2").Value = $statusText
$ws.Range("C2").Value = $statusText
$ws.Range("D2").Value = $handoffDate

$ws.Range("A3").Value = "$uuid2.md"
$ws.Range("B3").Value = $statusText
$ws.Range("C3").Value = $statusText
$ws.Range("D3").Value = $handoffDate

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null

# ---------- zh-cn sheet ----------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$uuid1.md"
$ws.Range("C2").Value = $statusText
$ws.Range("D2").Value = $xlf1zh
$ws.Range("E2").Value = $zhHandoffDatetime
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()
$ws.Range("H2").Value = $neverBack

$ws.Range("A3").Value = "$uuid2.md"
$ws.Range("C3").Value = $statusText
$ws.Range("D3").Value = $xlf1zh
$ws.Range("E3").Value = $zhHandoffDatetime
$ws.Range("F3").Clear()
$ws.Range("G3").Clear()
$ws.Range("H3").Value = $neverBack

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf1zh", "", "", $xlf1zh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf1zh", "", "", $xlf1zh) | Out-Null

# ---------- de-de sheet ----------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$uuid1.md"
$ws.Range("C2").Value = $statusText
$ws.Range("D2").Value = $xlf1de
$ws.Range("E2").Value = $deHandoffDatetime
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()
$ws.Range("H2").Value = $neverBack

$ws.Range("A3").Value = "$uuid2.md"
$ws.Range("C3").Value = $statusText
$ws.Range("D3").Value = $xlf1de
$ws.Range("E3").Value = $deHandoffDatetime
$ws.Range("F3").Clear()
$ws.Range("G3").Clear()
$ws.Range("H3").Value = $neverBack

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf1de", "", "", $xlf1de) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf1de", "", "", $xlf1de) | Out-Null
